# New Horizons Cinemas.docx -- "MySQL workbench file name change"
#
# The document already ends with the "Use a CASE tool, Computer Aided
# Software Engineering." paragraph. Two new paragraphs are appended
# after it, referencing the separate MySQL Workbench file that the
# ERD was built in:
#
#   "Please see the separate file. "
#   "New Horizons Cinemas MySQL Workbench.mwb"
#
# Both paragraphs use a 0.5in (720 twip / 36pt) left indent and
# complex-script bold (bCs) run formatting, matching the look of the
# other "reference" paragraphs already in the document (e.g. the
# existing "MySQL Workbench." line), with the bullet numbering that a
# freshly-split final paragraph would otherwise inherit stripped back
# off.

$d = $word.ActiveDocument

# --- Paragraph 1: "Please see the separate file. " -----------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$endOfDoc = $lastPara.Range
$endOfDoc.Collapse(0)
$endOfDoc.InsertParagraphAfter()

$p1 = $d.Paragraphs($d.Paragraphs.Count)
$p1.Style = "Normal"
$p1.Range.ParagraphFormat.LeftIndent = 36

$p1.Range.Text = "Please see the separate file. "

$p1 = $d.Paragraphs($d.Paragraphs.Count)
$t1 = $p1.Range
[void]$t1.MoveEnd(1, -1)
$t1.Font.Bold = 0
$t1.Font.BoldBi = 1

# --- Paragraph 2: "New Horizons Cinemas MySQL Workbench.mwb" -------
$afterP1 = $p1.Range
$afterP1.Collapse(0)
$afterP1.InsertParagraphAfter()

$p2 = $d.Paragraphs($d.Paragraphs.Count)
$p2.Style = "Normal"
$p2.Range.ParagraphFormat.LeftIndent = 36

$p2.Range.Text = "New Horizons Cinemas MySQL Workbench.mwb"

$p2 = $d.Paragraphs($d.Paragraphs.Count)
$t2 = $p2.Range
[void]$t2.MoveEnd(1, -1)
$t2.Font.Bold = 0
$t2.Font.BoldBi = 1

Write-Output "Appended MySQL Workbench filename paragraphs."
